$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 17 new rows above the current row 1, pushing all existing data down.
$ws.Rows("1:17").Insert()

# New dates for rows 1..17 (the 1st of each month, continuing the existing
# monthly series backwards from 2022-01-01).
$dates = @(44044,44075,44105,44136,44166,44197,44228,44256,44287,44317,44348,44378,44409,44440,44470,44501,44531)

# Copy the date-formatted style from the (now shifted) original column-A
# cells so the new cells inherit the same number format (mm/dd/yyyy style).
$ws.Range("A18:A34").Copy($ws.Range("A1:A17"))

for ($i = 1; $i -le 17; $i++) {
    $ws.Cells.Item($i, 1).Value = $dates[$i - 1]
}

# New column-B formulas for rows 1..17: each references the cell 17 rows
# below it (the original row before the insert), mirroring the alternating
# +/-200000 pattern used elsewhere in the sheet.
for ($i = 1; $i -le 17; $i++) {
    $r = $i + 17
    $ws.Cells.Item($i, 2).Formula = "=B" + $r + "+(-1)^(B" + $r + "/100000)*200000"
}

# Match the new selection left behind in the saved workbook.
$ws.Range("B1:B17").Select() | Out-Null

